$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.064621547438719
$ws.Range("C2").Value = 0.05893764867289519
$ws.Range("D2").Value = 0.005120677384145012
$ws.Range("F2").Value = 4.817453757367403
$ws.Range("G2").Value = 0.002663311908646904
$ws.Range("I2").Value = 3.169288271331212
$ws.Range("J2").Value = 0.1962860345525819
$ws.Range("K2").Value = 0.9242871356489104
$ws.Range("L2").Value = 0.3500053865620742
$ws.Range("M2").Value = 0.3086840539010964
$ws.Range("B3").Value = 1.051262052412625
$ws.Range("C3").Value = 0.0536964426258919
$ws.Range("D3").Value = 0.005181647594255878
$ws.Range("F3").Value = 4.776777517775415
$ws.Range("G3").Value = 0.002667603283948381
$ws.Range("I3").Value = 3.143608035752749
$ws.Range("J3").Value = 0.196287625653369
$ws.Range("K3").Value = 0.9046665186712062
$ws.Range("L3").Value = 0.349776915524771
$ws.Range("M3").Value = 0.3066444472197141
$ws.Range("B4").Value = 1.043799320071685
$ws.Range("C4").Value = 0.05051067835978529
$ws.Range("D4").Value = 0.005225134436640744
$ws.Range("F4").Value = 4.752988635236477
$ws.Range("G4").Value = 0.002670378449959713
$ws.Range("I4").Value = 3.128492370065374
$ws.Range("J4").Value = 0.1963080200593375
$ws.Range("K4").Value = 0.8932858180698702
$ws.Range("L4").Value = 0.349783151804175
$ws.Range("M4").Value = 0.3055743500323977
$ws.Range("B5").Value = 1.040944534901087
$ws.Range("C5").Value = 0.04922051551383788
$ws.Range("D5").Value = 0.005244380404027638
$ws.Range("F5").Value = 4.743592717174806
$ws.Range("G5").Value = 0.002671544730757571
$ws.Range("I5").Value = 3.122496450726274
$ws.Range("J5").Value = 0.1963212194635382
$ws.Range("K5").Value = 0.888815760176854
$ws.Range("L5").Value = 0.3498226076385933
$ws.Range("M5").Value = 0.3051841729131439
$ws.Range("B6").Value = 1.040481761805239
$ws.Range("C6").Value = 0.04900676971384144
$ws.Range("D6").Value = 0.005247668366351377
$ws.Range("F6").Value = 4.742050542365419
$ws.Range("G6").Value = 0.002671740531014659
$ws.Range("I6").Value = 3.121510724242086
$ws.Range("J6").Value = 0.1963237066414116
$ws.Range("K6").Value = 0.8880836400975198
$ws.Range("L6").Value = 0.3498313909054929
$ws.Range("M6").Value = 0.3051221584922814
$ws.Range("B7").Value = 1.043760064631925
$ws.Range("C7").Value = 0.05049324622797258
$ws.Range("D7").Value = 0.005225387816897964
$ws.Range("F7").Value = 4.752860711180361
$ws.Range("G7").Value = 0.002670394035303565
$ws.Range("I7").Value = 3.128410843828092
$ws.Range("J7").Value = 0.196308178270332
$ws.Range("K7").Value = 0.893224854295454
$ws.Range("L7").Value = 0.3497835343587923
$ws.Range("M7").Value = 0.3055689020379688
$ws.Range("B8").Value = 1.05986173421249
$ws.Range("C8").Value = 0.05712371835302577
$ws.Range("D8").Value = 0.005140445927821702
$ws.Range("F8").Value = 4.803182248659155
$ws.Range("G8").Value = 0.002664762535021457
$ws.Range("I8").Value = 3.16029828703978
$ws.Range("J8").Value = 0.196282556461469
$ws.Range("K8").Value = 0.9173836699124251
$ws.Range("L8").Value = 0.349896241503906
$ws.Range("M8").Value = 0.3079430176436659
$ws.Range("B9").Value = 1.097301926579121
$ws.Range("C9").Value = 0.07038687010046374
$ws.Range("D9").Value = 0.005021755935569416
$ws.Range("F9").Value = 4.911289945927336
$ws.Range("G9").Value = 0.00265482680055816
$ws.Range("I9").Value = 3.228014675919937
$ws.Range("J9").Value = 0.1963861720743445
$ws.Range("K9").Value = 0.9700474021820185
$ws.Range("L9").Value = 0.3512770982223969
$ws.Range("M9").Value = 0.3140421002911609
$ws.Range("B10").Value = 1.128380596269096
$ws.Range("C10").Value = 0.08029697193255458
$ws.Range("D10").Value = 0.004963582726777105
$ws.Range("F10").Value = 4.996493584983284
$ws.Range("G10").Value = 0.00264819502988069
$ws.Range("I10").Value = 3.28094916980848
$ws.Range("J10").Value = 0.1965558853105094
$ws.Range("K10").Value = 1.011970544735163
$ws.Range("L10").Value = 0.3529956346380914
$ws.Range("M10").Value = 0.3194008647576823
$ws.Range("B11").Value = 1.143294189931908
$ws.Range("C11").Value = 0.08484281030388274
$ws.Range("D11").Value = 0.00494338933159888
$ws.Range("F11").Value = 5.036516520812569
$ws.Range("G11").Value = 0.002645321580581728
$ws.Range("I11").Value = 3.305726881423496
$ws.Range("J11").Value = 0.1966533746423735
$ws.Range("K11").Value = 1.031746122981048
$ws.Range("L11").Value = 0.3539297215195205
$ws.Range("M11").Value = 0.3220289343403948
$ws.Range("B12").Value = 1.149053001950534
$ws.Range("C12").Value = 0.08656972017254816
$ws.Range("D12").Value = 0.004936641518455076
$ws.Range("F12").Value = 5.051854180877115
$ws.Range("G12").Value = 0.002644253980165576
$ws.Range("I12").Value = 3.315210140863584
$ws.Range("J12").Value = 0.1966932039405656
$ws.Range("K12").Value = 1.03933597964047
$ws.Range("L12").Value = 0.3543052804613609
$ws.Range("M12").Value = 0.323051437090939
$ws.Range("B13").Value = 1.147807789102558
$ws.Range("C13").Value = 0.08619755334817114
$ws.Range("D13").Value = 0.004938054839785266
$ws.Range("F13").Value = 5.048542851069953
$ws.Range("G13").Value = 0.002644482996561778
$ws.Range("I13").Value = 3.313163279736358
$ws.Range("J13").Value = 0.1966844965448935
$ws.Range("K13").Value = 1.037696864513293
$ws.Range("L13").Value = 0.3542234267450084
$ws.Range("M13").Value = 0.3228300091718523
$ws.Range("B14").Value = 1.143765740177116
$ws.Range("C14").Value = 0.0849847735337903
$ws.Range("D14").Value = 0.004942816183394427
$ws.Range("F14").Value = 5.037774713940934
$ws.Range("G14").Value = 0.002645233337867573
$ws.Range("I14").Value = 3.306505060172967
$ws.Range("J14").Value = 0.1966565930860931
$ws.Range("K14").Value = 1.032368515552548
$ws.Range("L14").Value = 0.3539601816176656
$ws.Range("M14").Value = 0.3221125092713706
$ws.Range("B15").Value = 1.141304364909701
$ws.Range("C15").Value = 0.08424262948955175
$ws.Range("D15").Value = 0.004945849641700306
$ws.Range("F15").Value = 5.031202605483799
$ws.Range("G15").Value = 0.002645695611876192
$ws.Range("I15").Value = 3.302439798729722
$ws.Range("J15").Value = 0.1966398805215874
$ws.Range("K15").Value = 1.029117935905731
$ws.Range("L15").Value = 0.3538017788199284
$ws.Range("M15").Value = 0.3216765746270553
$ws.Range("B16").Value = 1.127421551188917
$ws.Range("C16").Value = 0.08000065587197014
$ws.Range("D16").Value = 0.004965028354169121
$ws.Range("F16").Value = 4.993903428859682
$ws.Range("G16").Value = 0.002648385692649272
$ws.Range("I16").Value = 3.279343939115108
$ws.Range("J16").Value = 0.1965499219399227
$ws.Range("K16").Value = 1.010692337897865
$ws.Range("L16").Value = 0.3529376495019605
$ws.Range("M16").Value = 0.319232939132128
$ws.Range("B17").Value = 1.119103442816993
$ws.Range("C17").Value = 0.07740805125926897
$ws.Range("D17").Value = 0.004978397860025296
$ws.Range("F17").Value = 4.971345293056203
$ws.Range("G17").Value = 0.00265007261773048
$ws.Range("I17").Value = 3.265354157069495
$ws.Range("J17").Value = 0.1964999274855259
$ws.Range("K17").Value = 0.9995692496628124
$ws.Range("L17").Value = 0.3524465013001432
$ws.Range("M17").Value = 0.3177825549459605
$ws.Range("B18").Value = 1.114392115445582
$ws.Range("C18").Value = 0.07592039572530496
$ws.Range("D18").Value = 0.004986677950098439
$ws.Range("F18").Value = 4.958489354093871
$ws.Range("G18").Value = 0.00265105639396245
$ws.Range("I18").Value = 3.257373250287031
$ws.Range("J18").Value = 0.1964730816965847
$ws.Range("K18").Value = 0.9932378590098949
$ws.Range("L18").Value = 0.3521783462244201
$ws.Range("M18").Value = 0.3169662516505163
$ws.Range("B19").Value = 1.112809489701107
$ws.Range("C19").Value = 0.07541730702783411
$ws.Range("D19").Value = 0.004989582912649837
$ws.Range("F19").Value = 4.954156970814978
$ws.Range("G19").Value = 0.002651391805678396
$ws.Range("I19").Value = 3.254682321921521
$ws.Range("J19").Value = 0.1964643203158403
$ws.Range("K19").Value = 0.9911055498861003
$ws.Range("L19").Value = 0.3520900184111042
$ws.Range("M19").Value = 0.316692944922174
$ws.Range("B20").Value = 1.119981362897789
$ws.Range("C20").Value = 0.07768367127363263
$ws.Range("D20").Value = 0.004976913578607522
$ws.Range("F20").Value = 4.973734337744986
$ws.Range("G20").Value = 0.002649891645204724
$ws.Range("I20").Value = 3.266836597900735
$ws.Range("J20").Value = 0.1965050518701723
$ws.Range("K20").Value = 1.000746458270356
$ws.Range("L20").Value = 0.3524973011445667
$ws.Range("M20").Value = 0.31793509668098
$ws.Range("B21").Value = 1.144949967225955
$ws.Range("C21").Value = 0.08534084646129259
$ws.Range("D21").Value = 0.00494139328307952
$ws.Range("F21").Value = 5.040932640290009
$ws.Range("G21").Value = 0.002645012388506438
$ws.Range("I21").Value = 3.308458012342015
$ws.Range("J21").Value = 0.1966647100108183
$ws.Range("K21").Value = 1.033930832503557
$ws.Range("L21").Value = 0.3540369108389569
$ws.Range("M21").Value = 0.3223225156822096
$ws.Range("B22").Value = 1.161917386223962
$ws.Range("C22").Value = 0.09037734658011232
$ws.Range("D22").Value = 0.004923417538565289
$ws.Range("F22").Value = 5.085910617346826
$ws.Range("G22").Value = 0.002641943029027348
$ws.Range("I22").Value = 3.336245754895671
$ws.Range("J22").Value = 0.1967860281466791
$ws.Range("K22").Value = 1.056209031621279
$ws.Range("L22").Value = 0.3551704097077391
$ws.Range("M22").Value = 0.3253491029637914
$ws.Range("B23").Value = 1.152802235034073
$ws.Range("C23").Value = 0.08768630825883861
$ws.Range("D23").Value = 0.004932533026231312
$ws.Range("F23").Value = 5.061807976684946
$ws.Range("G23").Value = 0.002643570302289626
$ws.Range("I23").Value = 3.321361261220261
$ws.Range("J23").Value = 0.1967197270110148
$ws.Range("K23").Value = 1.044264736013702
$ws.Range("L23").Value = 0.3545538147302807
$ws.Range("M23").Value = 0.3237192134196221
$ws.Range("B24").Value = 1.119584234477145
$ws.Range("C24").Value = 0.07755905454112622
$ws.Range("D24").Value = 0.00497758277257887
$ws.Range("F24").Value = 4.972653898528222
$ws.Range("G24").Value = 0.002649973419388276
$ws.Range("I24").Value = 3.266166193248438
$ws.Range("J24").Value = 0.1965027292280439
$ws.Range("K24").Value = 1.000214044704848
$ws.Range("L24").Value = 0.3524742902595932
$ws.Range("M24").Value = 0.3178660779162001
$ws.Range("B25").Value = 1.086545927406775
$ws.Range("C25").Value = 0.06677022777914488
$ws.Range("D25").Value = 0.005120677384145012
$ws.Range("F25").Value = 4.881031626878297
$ws.Range("G25").Value = 0.002657396850617915
$ws.Range("I25").Value = 3.209138193912452
$ws.Range("J25").Value = 0.1963416914170537
$ws.Range("K25").Value = 0.9552335936506893
$ws.Range("L25").Value = 0.3507796895964077
$ws.Range("M25").Value = 0.3122378379117627

Write-Output "Applied 240 cell updates"
